$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the wrong deployed date in cell A1
$ws.Range("A1").Value = "Date Deployed: 10/3/2019"

# Update the selected cell to A2
$ws.Range("A2").Select()
